# Updates the crypto price/volume figures on Sheet1 to the refreshed
# values pulled on Tue Feb  7 20:53:18 UTC 2023 (GitHub Actions symbol
# list refresh). Only the "Price" (D) and "Volume(1h)" (E) columns for
# a handful of rows change; everything else (labels, links, dates,
# styles) stays untouched.
#
# Cells in this sheet are stored as plain text (e.g. "0.52%", "44.43"),
# not real numbers/percentages. Assigning a string that merely *looks*
# like a number or percentage straight to Range.Value/.Value2/.Formula
# makes Excel "smart" parse it into a real number and silently stamp a
# new number-format style onto the cell - which would corrupt the
# original (style-less) formatting of these cells. To avoid that, each
# new value is entered as a text-producing formula (="...") and then
# immediately collapsed to a literal value via Copy + PasteSpecial
# (values only), which preserves the text representation and leaves
# cell styling completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2"  = "0.52%"
    "D3"  = "44.43"
    "E3"  = "0.67%"
    "D4"  = "5.538"
    "E4"  = "-0.73%"
    "D5"  = "0.08163"
    "E5"  = "0.93%"
    "D6"  = "2.065"
    "E6"  = "4.40%"
    "D7"  = "0.9744"
    "E7"  = "2.31%"
    "E8"  = "-4.27%"
    "E9"  = "1.76%"
    "D10" = "10.22"
    "E10" = "-18.60%"
    "D11" = "0.1004"
    "E11" = "1.34%"
    "D12" = "0.04753"
    "E12" = "0.18%"
    "D13" = "0.1059"
    "E13" = "-0.85%"
    "D14" = "0.001265"
    "E14" = "-1.64%"
    "D15" = "0.04094"
    "E15" = "-3.19%"
    "D16" = "0.006013"
    "E16" = "0.40%"
    "D17" = "3.348"
    "E17" = "-0.77%"
    "E18" = "2.40%"
    "D19" = "2.645"
    "E19" = "2.96%"
    "D20" = "0.3350"
    "E20" = "-3.49%"
    "D21" = "0.1390"
    "E21" = "-1.40%"
    "D22" = "0.2567"
    "E22" = "2.36%"
    "E23" = "3.71%"
    "D24" = "0.004401"
    "E24" = "1.13%"
    "E25" = "7.15%"
    "E26" = "-6.23%"
    "D38" = "0.02681"
    "E38" = "0.83%"
    "D39" = "0.05659"
    "E39" = "2.14%"
    "E40" = "0.48%"
    "D41" = "0.1420"
    "E41" = "0.91%"
    "D42" = "0.007516"
    "E42" = "-7.09%"
    "D43" = "0.001953"
    "E43" = "-3.20%"
    "D44" = "0.008276"
    "E44" = "-6.99%"
    "D45" = "0.00007038"
    "E45" = "-2.98%"
    "E46" = "-0.38%"
    "E48" = "10.63%"
    "D49" = "0.003533"
    "E49" = "-26.17%"
    "D50" = "0.00002094"
    "E50" = "-0.38%"
    "E51" = "-0.38%"
}

foreach ($ref in $updates.Keys) {
    $newValue = $updates[$ref]
    $escaped = $newValue.Replace('"', '""')
    $cell = $ws.Range($ref)

    # Write the target text via a formula that evaluates to that exact
    # string, so Excel's autodetection cannot reinterpret it as a
    # number/percentage/date.
    $cell.Formula = '="' + $escaped + '"'

    # Freeze the formula result into a literal value (keeps the text
    # representation, drops the formula, and does not touch styling).
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = $false
